# Update the "想去人数" (F column) figures to the freshly generated output.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1651
$ws1.Range("F3").Value = 9064
$ws1.Range("F7").Value = 1076
$ws1.Range("F8").Value = 192
$ws1.Range("F11").Value = 5764
$ws1.Range("F15").Value = 4371
$ws1.Range("F17").Value = 160
$ws1.Range("F21").Value = 18
$ws1.Range("F24").Value = 2709

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1651
$ws4.Range("F3").Value = 9064
$ws4.Range("F5").Value = 16
$ws4.Range("F8").Value = 1076
$ws4.Range("F9").Value = 192
$ws4.Range("F12").Value = 5764
$ws4.Range("F16").Value = 4371
$ws4.Range("F18").Value = 160
$ws4.Range("F22").Value = 18
$ws4.Range("F25").Value = 2709
